$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("StatQuery"). This shifts the
# existing "dbExcel"/Neo4jData column from B to C and the existing
# "WebExcel"/WebData column from C to D, preserving their values,
# styles and widths.
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column B.
$ws.Range("B1").Value = "StatQuery"

# Query text for the newly inserted column B, row 2 (Insert() already
# copied the wrap-text style from column A onto the new B2 cell, same
# as Excel does when inserting a column in the middle of formatted
# data).
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Adenocarcinoma of the cervix'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Give the new column the same width as column A.
$ws.Columns.Item(2).ColumnWidth = 75

# Move/collapse the selection to A4, matching the saved workbook state.
$ws.Range("A4").Select() | Out-Null
